$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 362.16666
$ws.Range("I2").Value = 344.6
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 344.6
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -231.6
$ws.Range("N2").Value = -676

$ws.Range("H40").Value = 2472.5454
$ws.Range("I40").Value = 2174.75
$ws.Range("J40").Value = 2642.7144
$ws.Range("K40").Value = 2174.75
$ws.Range("L40").Value = 2642.7144
$ws.Range("M40").Value = -1999.75
$ws.Range("N40").Value = -2992.7144

$ws.Range("H74").Value = 3946.4285
$ws.Range("I74").Value = 2875
$ws.Range("J74").Value = 4750
$ws.Range("K74").Value = 2875
$ws.Range("L74").Value = 4750
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -6622

$ws.Range("H77").Value = 3946.4285
$ws.Range("I77").Value = 2875
$ws.Range("J77").Value = 4750
$ws.Range("K77").Value = 14375
$ws.Range("L77").Value = 23750
$ws.Range("M77").Value = -9695
$ws.Range("N77").Value = -33110

$ws.Range("H98").Value = 838.8
$ws.Range("I98").Value = 838.8
$ws.Range("K98").Value = 838.8
$ws.Range("M98").Value = 659.2

$ws.Range("H116").Value = 26566734
$ws.Range("I116").Value = 22820574
$ws.Range("K116").Value = 22820574
$ws.Range("M116").Value = -22817132

$ws.Range("H122").Value = 838.8
$ws.Range("I122").Value = 838.8
$ws.Range("K122").Value = 2516.4
$ws.Range("M122").Value = -66.39999999999964

$ws.Range("H129").Value = 1575.5555
$ws.Range("J129").Value = 1568.3334
$ws.Range("L129").Value = 4705.0002
$ws.Range("N129").Value = -14705.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1327
$ws.Range("I2").Value = 1196.6666
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 1196.6666
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -1083.6666
$ws.Range("N2").Value = -2726

$ws.Range("H74").Value = 1875.875
$ws.Range("I74").Value = 1215.3572
$ws.Range("J74").Value = 6499.5
$ws.Range("K74").Value = 1215.3572
$ws.Range("L74").Value = 6499.5
$ws.Range("M74").Value = -341.3571999999999
$ws.Range("N74").Value = -8247.5

$ws.Range("H77").Value = 1875.875
$ws.Range("I77").Value = 1215.3572
$ws.Range("J77").Value = 6499.5
$ws.Range("K77").Value = 6076.786
$ws.Range("L77").Value = 32497.5
$ws.Range("M77").Value = -1708.786
$ws.Range("N77").Value = -41233.5

$ws.Range("H116").Value = 1327
$ws.Range("I116").Value = 1196.6666
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 1196.6666
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 1097.3334
$ws.Range("N116").Value = -7088

$ws.Range("H132").Value = 37039470
$ws.Range("I132").Value = 50001964
$ws.Range("J132").Value = 3771.4285
$ws.Range("K132").Value = 150005892
$ws.Range("L132").Value = 11314.2855
$ws.Range("M132").Value = -150003362
$ws.Range("N132").Value = -16374.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1327
$ws.Range("I3").Value = 1196.6666
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 1196.6666
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -1082.6666
$ws.Range("N3").Value = -2728

$ws.Range("H64").Value = 13498.375
$ws.Range("I64").Value = 646.3333
$ws.Range("J64").Value = 21209.6
$ws.Range("K64").Value = 646.3333
$ws.Range("L64").Value = 21209.6
$ws.Range("M64").Value = -421.3333
$ws.Range("N64").Value = -21659.6

$ws.Range("H67").Value = 13498.375
$ws.Range("I67").Value = 646.3333
$ws.Range("J67").Value = 21209.6
$ws.Range("K67").Value = 646.3333
$ws.Range("L67").Value = 21209.6
$ws.Range("M67").Value = 133.6667
$ws.Range("N67").Value = -22769.6

$ws.Range("H82").Value = 49829.668
$ws.Range("J82").Value = 63994.5
$ws.Range("L82").Value = 63994.5
$ws.Range("N82").Value = -64760.5

$ws.Range("H85").Value = 49829.668
$ws.Range("J85").Value = 63994.5
$ws.Range("L85").Value = 63994.5
$ws.Range("N85").Value = -66646.5

$ws.Range("H86").Value = 27780748
$ws.Range("I86").Value = 41670016
$ws.Range("K86").Value = 41670016
$ws.Range("M86").Value = -41668893

$ws.Range("H89").Value = 27780748
$ws.Range("I89").Value = 41670016
$ws.Range("K89").Value = 208350080
$ws.Range("M89").Value = -208344464

$ws.Range("H105").Value = 1984.4286
$ws.Range("I105").Value = 1932.5333
$ws.Range("K105").Value = 1932.5333
$ws.Range("M105").Value = -185.5333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3086.4707
$ws.Range("I132").Value = 2778.5
$ws.Range("K132").Value = 8335.5
$ws.Range("M132").Value = -5805.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 446.3889
$ws.Range("J122").Value = 572.375
$ws.Range("L122").Value = 5151.375
$ws.Range("N122").Value = -10051.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 138
$ws.Range("I2").Value = 64.8
$ws.Range("J2").Value = 199
$ws.Range("K2").Value = 64.8
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = 48.2
$ws.Range("N2").Value = -425

$ws.Range("H70").Value = 12345.454
$ws.Range("I70").Value = 11756.889
$ws.Range("K70").Value = 11756.889
$ws.Range("M70").Value = -11486.889

$ws.Range("H73").Value = 12345.454
$ws.Range("I73").Value = 11756.889
$ws.Range("K73").Value = 11756.889
$ws.Range("M73").Value = -10820.889

$ws.Range("H80").Value = 3200.4
$ws.Range("J80").Value = 3166.5
$ws.Range("L80").Value = 3166.5
$ws.Range("N80").Value = -5162.5

$ws.Range("H83").Value = 3200.4
$ws.Range("J83").Value = 3166.5
$ws.Range("L83").Value = 15832.5
$ws.Range("N83").Value = -25816.5

$ws.Range("H113").Value = 1216.2858
$ws.Range("I113").Value = 1079.8462
$ws.Range("K113").Value = 1079.8462
$ws.Range("M113").Value = 1090.1538

$ws.Range("H122").Value = 2430.2778
$ws.Range("I122").Value = 2117
$ws.Range("J122").Value = 3996.6667
$ws.Range("K122").Value = 6351
$ws.Range("L122").Value = 11990.0001
$ws.Range("M122").Value = -3901
$ws.Range("N122").Value = -16890.0001

$ws.Range("H126").Value = 17126.5
$ws.Range("I126").Value = 21502
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 64506
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -62036
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1780.7693
$ws.Range("I16").Value = 1660.7778
$ws.Range("K16").Value = 1660.7778
$ws.Range("M16").Value = -1490.7778

$ws.Range("H46").Value = 1865.95
$ws.Range("I46").Value = 992.625
$ws.Range("J46").Value = 2448.1667
$ws.Range("K46").Value = 992.625
$ws.Range("L46").Value = 2448.1667
$ws.Range("M46").Value = -804.625
$ws.Range("N46").Value = -2824.1667

$ws.Range("H61").Value = 21134.215
$ws.Range("I61").Value = 18402.6
$ws.Range("J61").Value = 22651.777
$ws.Range("K61").Value = 18402.6
$ws.Range("L61").Value = 22651.777
$ws.Range("M61").Value = -18200.6
$ws.Range("N61").Value = -23055.777

$ws.Range("H93").Value = 1378.1111
$ws.Range("I93").Value = 1300.4286
$ws.Range("K93").Value = 1300.4286
$ws.Range("M93").Value = -52.42859999999996

$ws.Range("H100").Value = 3627.8572
$ws.Range("I100").Value = 3498.889
$ws.Range("J100").Value = 3860
$ws.Range("K100").Value = 3498.889
$ws.Range("L100").Value = 3860
$ws.Range("M100").Value = -2957.889
$ws.Range("N100").Value = -4942

$ws.Range("H113").Value = 21134.215
$ws.Range("I113").Value = 18402.6
$ws.Range("J113").Value = 22651.777
$ws.Range("K113").Value = 18402.6
$ws.Range("L113").Value = 22651.777
$ws.Range("M113").Value = -16232.6
$ws.Range("N113").Value = -26991.777

$ws.Range("H132").Value = 3466.4583
$ws.Range("I132").Value = 3170
$ws.Range("K132").Value = 9510
$ws.Range("M132").Value = -6980

$ws.Range("H136").Value = 2359
$ws.Range("I136").Value = 2190.8
$ws.Range("K136").Value = 6572.400000000001
$ws.Range("M136").Value = -4022.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 468.66666
$ws.Range("J107").Value = 591
$ws.Range("L107").Value = 1773
$ws.Range("N107").Value = -5613

$ws.Range("H136").Value = 7111.5
$ws.Range("I136").Value = 5500
$ws.Range("J136").Value = 8262.571
$ws.Range("K136").Value = 16500
$ws.Range("L136").Value = 24787.713
$ws.Range("M136").Value = -13950
$ws.Range("N136").Value = -29887.713
